$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: remove all existing content/formatting on the sheet.
$ws.Cells.Clear()

function Set-PlainCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = $val
}

function Set-SmallCell($addr, $val, $fmt) {
    $c = $ws.Range($addr)
    $c.Value = $val
    $c.Font.Name = "Arial"
    $c.Font.Size = 9
    if ($fmt) {
        $c.NumberFormat = $fmt
    }
}

# ---- Row 1: headers -------------------------------------------------------
# A1:E1 keep the workbook's default (unstyled) look.
Set-PlainCell "A1" "idx"
Set-PlainCell "B1" "idx2"
Set-PlainCell "C1" "Name"
Set-PlainCell "D1" "Date Start"
Set-PlainCell "E1" "Date End"

# F1:K1 use the small Arial (9pt) font used throughout the rest of the table.
Set-SmallCell "F1" "(m3/s)" $null
Set-SmallCell "G1" "(MW1)" $null
Set-SmallCell "H1" "(MW2)" $null
Set-SmallCell "I1" "(GWh) Winter" $null
Set-SmallCell "J1" "(GWh) Summer" $null
Set-SmallCell "K1" "(GWh) Year" $null

# ---- Data rows -------------------------------------------------------------
# Columns A, B, D, E -> integer number format "0"
# Column C             -> plain Arial 9pt text (General format)
# Columns F..K          -> 2-decimal number format "0.00"

# Row 2 - Wunderklingen
Set-SmallCell "A2" 1 "0"
Set-SmallCell "B2" 108900 "0"
Set-SmallCell "C2" "Wunderklingen" $null
Set-SmallCell "D2" 1895 "0"
Set-SmallCell "E2" 1968 "0"
Set-SmallCell "F2" 5.5 "0.00"
Set-SmallCell "G2" 0.42 "0.00"
Set-SmallCell "H2" 0.41 "0.00"
Set-SmallCell "I2" 1.4 "0.00"
Set-SmallCell "J2" 1 "0.00"
Set-SmallCell "K2" 2.4 "0.00"

# Row 3 - Engeweiher
Set-SmallCell "A3" 2 "0"
Set-SmallCell "B3" 106300 "0"
Set-SmallCell "C3" "Engeweiher" $null
Set-SmallCell "D3" 1909 "0"
Set-SmallCell "E3" 1993 "0"
Set-SmallCell "F3" 4 "0.00"
Set-SmallCell "G3" 5 "0.00"
Set-SmallCell "H3" 5 "0.00"

# Row 4 - Eglisau
Set-SmallCell "A4" 3 "0"
Set-SmallCell "B4" 108700 "0"
Set-SmallCell "C4" "Eglisau" $null
Set-SmallCell "D4" 1920 "0"
Set-SmallCell "E4" 1927 "0"
Set-SmallCell "F4" 400 "0.00"
Set-SmallCell "G4" 10.82 "0.00"
Set-SmallCell "H4" 10.34 "0.00"
Set-SmallCell "I4" 37.49 "0.00"
Set-SmallCell "J4" 39.59 "0.00"
Set-SmallCell "K4" 77.08 "0.00"

# Row 5 - Neuhausen (no "Date End")
Set-SmallCell "A5" 4 "0"
Set-SmallCell "B5" 106400 "0"
Set-SmallCell "C5" "Neuhausen" $null
Set-SmallCell "D5" 1951 "0"
Set-SmallCell "F5" 25 "0.00"
Set-SmallCell "G5" 2.32 "0.00"
Set-SmallCell "H5" 2.2 "0.00"
Set-SmallCell "I5" 9.95 "0.00"
Set-SmallCell "J5" 9.95 "0.00"
Set-SmallCell "K5" 19.9 "0.00"

# Row 6 - Rheinau
Set-SmallCell "A6" 5 "0"
Set-SmallCell "B6" 106500 "0"
Set-SmallCell "C6" "Rheinau" $null
Set-SmallCell "D6" 1956 "0"
Set-SmallCell "E6" 2005 "0"
Set-SmallCell "F6" 400 "0.00"
Set-SmallCell "G6" 2.98 "0.00"
Set-SmallCell "H6" 2.92 "0.00"
Set-SmallCell "I6" 6.39 "0.00"
Set-SmallCell "J6" 13.19 "0.00"
Set-SmallCell "K6" 19.58 "0.00"

# Row 7 - Schaffhausen (no "Date End")
Set-SmallCell "A7" 6 "0"
Set-SmallCell "B7" 106200 "0"
Set-SmallCell "C7" "Schaffhausen" $null
Set-SmallCell "D7" 1964 "0"
Set-SmallCell "F7" 500 "0.00"
Set-SmallCell "G7" 22.57 "0.00"
Set-SmallCell "H7" 19.84 "0.00"
Set-SmallCell "I7" 62.06 "0.00"
Set-SmallCell "J7" 73.64 "0.00"
Set-SmallCell "K7" 135.7 "0.00"

# ---- Sheet view: active cell / selection ----------------------------------
$ws.Range("A4:K4").Select()
